$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was date 44505 / Californiana(o) / Primera / 100 / $/bandeja) ->
#        date 44902 / Golden Nugget / Especial / 60 / $/caja
$ws.Range("D2").Value = 44902
$ws.Range("K2").Value = "Golden Nugget"
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 60
$ws.Range("Q2").Value = "$/caja 10 kilos"

# Row 3 (was date 44505 / Golden Nugget / Primera / 50 / 15000 / $/bandeja / 1500) ->
#        date 44902 / Golden Nugget / Primera / 70 / 13000 / $/caja / 1300
$ws.Range("D3").Value = 44902
$ws.Range("M3").Value = 70
$ws.Range("N3").Value = 13000
$ws.Range("O3").Value = 13000
$ws.Range("P3").Value = 13000
$ws.Range("Q3").Value = "$/caja 10 kilos"
$ws.Range("S3").Value = 1300

# Row 4 (was date 44902 / Golden Nugget / Especial / 60 / $/caja) ->
#        date 44505 / Californiana(o) / Primera / 100 / $/bandeja
$ws.Range("D4").Value = 44505
$ws.Range("K4").Value = "Californiana(o)"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("Q4").Value = "$/bandeja 10 kilos"

# Row 5 (was date 44902 / Golden Nugget / Primera / 70 / 13000 / $/caja / 1300) ->
#        date 44505 / Golden Nugget / Primera / 50 / 15000 / $/bandeja / 1500
$ws.Range("D5").Value = 44505
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = "$/bandeja 10 kilos"
$ws.Range("S5").Value = 1500
